$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.485.78"
$ws.Range("E2").Value = "  +3.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.51"
$ws.Range("E3").Value = "  +4.31%  "
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.91"
$ws.Range("E5").Value = "  +3.22%  "
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3840"
$ws.Range("E7").Value = "  +2.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3548"
$ws.Range("E8").Value = "  +3.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.93"
$ws.Range("E9").Value = "  -2.11%  "
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07798"
$ws.Range("E11").Value = "  +3.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.40"
$ws.Range("E13").Value = "  +9.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.614"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.817.82"
$ws.Range("E15").Value = "  +4.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.203"
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001125"
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06731"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.66"
$ws.Range("E21").Value = "  +5.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.568"
$ws.Range("E22").Value = "  +6.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.20"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.486.23"
$ws.Range("E24").Value = "  +3.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.469"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.699"
$ws.Range("E26").Value = "  +7.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.24"
$ws.Range("E27").Value = "  +13.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.471"
$ws.Range("E28").Value = "  +3.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "154.14"
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.019.96"
$ws.Range("E30").Value = "  +4.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "136.36"
$ws.Range("E31").Value = "  +3.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.379"
$ws.Range("E32").Value = "  +2.63%  "
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.97"
$ws.Range("E34").Value = "  +5.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08808"
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.689"
$ws.Range("E36").Value = "  -2.08%  "
$ws.Range("E37").Value = "  +2.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.7042"
$ws.Range("E38").Value = "  +12.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06530"
$ws.Range("E39").Value = "  +2.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2261"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02402"
$ws.Range("E41").Value = "  +1.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.002"
$ws.Range("E42").Value = "  +3.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.299"
$ws.Range("E43").Value = "  +4.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.96"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6622"
$ws.Range("E45").Value = "  +8.28%  "
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.964"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.198"
$ws.Range("E48").Value = "  +5.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.56"
$ws.Range("E49").Value = "  +2.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07328"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.97"
$ws.Range("E51").Value = "  +3.86%  "
